$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").NumberFormat = "DD/MM/YY"
$ws.Range("E2").Value = 42844
$ws.Range("F2").NumberFormat = "DD/MM/YY"
$ws.Range("F2").Value = 42845

$ws.Columns("E").ColumnWidth = 10.25

$ws.Range("F14").Select()
